# Add a new "plot_labels" worksheet after "values" containing GLORIA
# Enduse_name -> plot_label lookups, then update sheet selections/view state
# to match: "values" is no longer the active tab, "plot_labels" is.

$wb = $excel.ActiveWorkbook

$valuesSheet = $wb.Worksheets.Item("values")

$ws = $wb.Worksheets.Add($null, $valuesSheet)
$ws.Name = "plot_labels"

$colA = @('Enduse_name','Buildings','Civil_engineering_except_roads','Machinery_and_equipment','Computers_and_precision_instruments','Electrical_equipment','Motor_vehicles_trailers_and_semi.trailers','Other_transport_equipment','Furniture_and_other_manufactured_goods_nec','Textiles','Printed_matter_and_recorded_media','Food_packaging','Products_nec','Construction','Machinery and equipment n.e.c. ','Office machinery and computers','Radio, television and communication equipment and apparatus','Medical, precision and optical instruments, watches and clocks','Electrical machinery and apparatus n.e.c.','Motor vehicles, trailers and semi-trailers','Other transport equipment','Furniture; other manufactured goods n.e.c.','Textiles','Printed matter and recorded media','Food','Products nec','Other raw materials','Secondary materials','Energy carriers','Other','Services')
$colB = @('plot_label','Buildings','Civil engineering except roads','Machinery and equipment','Computers and precision instruments','Electrical equipment','Motor vehicles trailers and semi.trailers','Other transport equipment','Furniture and other manufactured goods nec','Textiles','Printed matter and recorded_media','Food packaging','Products nec','Construction','Machinery and equipment n.e.c. ','Office machinery and computers','Radio, television and communication equipment and apparatus','Medical, precision and optical instruments, watches and clocks','Electrical machinery and apparatus n.e.c.','Motor vehicles, trailers and semi-trailers','Other transport equipment','Furniture; other manufactured goods n.e.c.','Textiles','Printed matter and recorded media','Food','Products nec','Other raw materials','Secondary materials','Energy carriers','Other','Services')

for ($i = 0; $i -lt $colA.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $colA[$i]
    $ws.Cells.Item($i + 1, 2).Value = $colB[$i]
}

$ws.Columns.Item(1).ColumnWidth = 35.666666666666664
$ws.Columns.Item(2).ColumnWidth = 45.83

# restore the "values" sheet selection/view (it is no longer the active tab)
$valuesSheet.Activate()
$valuesSheet.Range("B2:B20").Select()

# make the new sheet the active / visible one, matching the new activeTab
$ws.Activate()
$ws.Range("B10").Select()

Write-Output "done"
